$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold plain text (prices/percentages with unusual
# formatting like "41.365.24" or "0.0₃0954"); force text format so Excel
# does not reinterpret the assigned strings as numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "41.365.24"
$ws.Range("E2").Value = "  -5.71%  "
$ws.Range("D3").Value = "2.219.56"
$ws.Range("E3").Value = "  -6.19%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "243.62"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  -7.92%  "
$ws.Range("D7").Value = "69.81"
$ws.Range("E7").Value = "  -6.22%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.559"
$ws.Range("E9").Value = "  -7.24%  "
$ws.Range("D10").Value = "39.05"
$ws.Range("E10").Value = "  +5.15%  "
$ws.Range("E11").Value = "  -7.15%  "
$ws.Range("D12").Value = "58.11"
$ws.Range("E12").Value = "  -3.47%  "
$ws.Range("E13").Value = "  -3.40%  "
$ws.Range("D14").Value = "6.74"
$ws.Range("E14").Value = "  -7.46%  "
$ws.Range("D15").Value = "2.549.93"
$ws.Range("E15").Value = "  -6.19%  "
$ws.Range("D16").Value = "14.83"
$ws.Range("E16").Value = "  -9.54%  "
$ws.Range("E17").Value = "  -9.62%  "
$ws.Range("D18").Value = "2.221.22"
$ws.Range("E18").Value = "  -6.24%  "
$ws.Range("D19").Value = "41.281.08"
$ws.Range("E19").Value = "  -5.83%  "
$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  -8.41%  "
$ws.Range("D21").Value = "72.21"
$ws.Range("E21").Value = "  -6.67%  "
$ws.Range("E22").Value = "  -7.98%  "
$ws.Range("D23").Value = "232.02"
$ws.Range("E23").Value = "  -7.78%  "
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  +10.33%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  -5.14%  "
$ws.Range("E27").Value = "  -3.32%  "
$ws.Range("D28").Value = "9.81"
$ws.Range("E28").Value = "  -7.40%  "
$ws.Range("E29").Value = "  -5.09%  "
$ws.Range("D30").Value = "172.44"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").Value = "20.47"
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  -8.17%  "
$ws.Range("D34").Value = "0.0717"
$ws.Range("E34").Value = "  -5.64%  "
$ws.Range("D35").Value = "5.22"
$ws.Range("E35").Value = "  -4.50%  "
$ws.Range("D36").Value = "4.60"
$ws.Range("E36").Value = "  -9.95%  "
$ws.Range("D37").Value = "3.92"
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("D38").Value = "24.15"
$ws.Range("E38").Value = "  +15.63%  "
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").Value = "2.28"
$ws.Range("E40").Value = "  -4.95%  "
$ws.Range("E41").Value = "  -11.53%  "
$ws.Range("D42").Value = "65.44"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").Value = "5.00"
$ws.Range("E43").Value = "  -11.04%  "
$ws.Range("D44").Value = "0.202"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "8.78"
$ws.Range("E45").Value = "  -3.41%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "10.92"
$ws.Range("E46").Value = "  +13.05%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.101"
$ws.Range("E47").Value = "  -6.39%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "4.58"
$ws.Range("E49").Value = "  +5.06%  "
$ws.Range("E50").Value = "  -6.04%  "
$ws.Range("E51").Value = "  -5.49%  "
